# Ticket 86 - Move, remove, and copy existing Excel hyperlinks.
# Adds a new "Shift" worksheet at the end of the workbook that demonstrates
# a jt:for loop copying/shifting a hyperlink ("Example" / "JETT") down the
# sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new "Shift" worksheet as the LAST sheet in the workbook -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Shift"

# --- Populate the cells ------------------------------------------------
$ws.Range("A1").Value = "Shift!"
$ws.Range("B1").Value = "Copy!"
$ws.Range("A2").Value = '<jt:for var="i" start="1" end="10">'
$ws.Range("B2").Value = "Example"
$ws.Range("C2").Value = "</jt:for>"
$ws.Range("B3").Value = "JETT"

# --- Add hyperlinks (B2 = "Example", B3 = "JETT") -----------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "http://sourceforge.net/projects/jett/")
$ws.Hyperlinks.Add($ws.Range("B3"), "http://sourceforge.net/projects/jett/")

# --- Restore the original active sheet/tab ------------------------------
# Adding a worksheet makes it the active sheet; the "Hyperlinks" sheet was
# (and remains) the one selected/shown when the workbook is opened.
$wb.Worksheets.Item(1).Activate()
